$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the sheet view / selection ---------------------------------
# Diff removes topLeftCell="B1" from the sheetView and changes the
# selection from H15 to B13.
$ws.Activate()
$ws.Range("B13").Select()

# --- Populate column H with 2020 data ----------------------------------
# For each source row, copy the formatting from an existing cell that
# already carries the desired number format/style into the new column H
# cell (so the style id matches the one produced by the original author),
# then set the value where applicable.
#
# Most rows simply reuse the style already present on column G of the
# same row. Rows 8 and 20, however, end up with the "0.0" numeric style
# (as used by column G on rows 10/22/23) rather than the General style
# that column G itself uses on those two rows, so those two are copied
# from a different source cell.
$rowsWithValues = [ordered]@{
    4  = @{ Value = 2020; StyleSource = "G4"  }
    5  = @{ Value = 42.2; StyleSource = "G5"  }
    7  = @{ Value = 42.5; StyleSource = "G7"  }
    8  = @{ Value = 42;   StyleSource = "G10" }
    10 = @{ Value = 50.9; StyleSource = "G10" }
    11 = @{ Value = 36.9; StyleSource = "G11" }
    12 = @{ Value = 34.8; StyleSource = "G12" }
    14 = @{ Value = 30.7; StyleSource = "G14" }
    15 = @{ Value = 48.8; StyleSource = "G15" }
    17 = @{ Value = 61.1; StyleSource = "G17" }
    18 = @{ Value = 56.7; StyleSource = "G18" }
    19 = @{ Value = 41.6; StyleSource = "G19" }
    20 = @{ Value = 49;   StyleSource = "G10" }
    21 = @{ Value = 43.5; StyleSource = "G21" }
    22 = @{ Value = 33.9; StyleSource = "G22" }
    23 = @{ Value = 34.6; StyleSource = "G23" }
    24 = @{ Value = 23.6; StyleSource = "G24" }
    25 = @{ Value = 35.9; StyleSource = "G25" }
}

$emptyRows = @(6, 9, 13, 16)

foreach ($r in $rowsWithValues.Keys) {
    $info = $rowsWithValues[$r]
    $srcCell = $ws.Range($info.StyleSource)
    $dstCell = $ws.Range("H$r")
    $srcCell.Copy()
    $dstCell.PasteSpecial(-4122)  # xlPasteFormats
    $dstCell.Value = $info.Value
}

foreach ($r in $emptyRows) {
    $srcCell = $ws.Range("G$r")
    $dstCell = $ws.Range("H$r")
    $srcCell.Copy()
    $dstCell.PasteSpecial(-4122)  # xlPasteFormats
}

$excel.CutCopyMode = 0

# Restore the originally requested selection after the paste operations.
$ws.Range("B13").Select()
